$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: new yarn-receive entries
$ws.Range("J8").Value = 44460
$ws.Range("J8").NumberFormat = "d-mmm"
$ws.Range("K8").Value = 1748
$ws.Range("L8").Value = "34/24"
$ws.Range("M8").Value = 72
$ws.Range("O8").Value = 1144

# Row 31 totals
$ws.Range("F31").Formula = "=SUM(F7:F30)"
$ws.Range("G31").Formula = "=SUM(G7:G30)"
$ws.Range("O31").Formula = "=SUM(O7:O30)"

# Row 7 short/excess formula
$ws.Range("P7").Formula = "=G31-O31"

# Scroll/selection state
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("P8").Select()
